$p = $ppt.ActivePresentation
$fonts = $p.Fonts
Write-Output $fonts.Count
for ($i = 1; $i -le $fonts.Count; $i++) {
    $f = $fonts.Item($i)
    Write-Output $f.Name
}
